$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff": rows for a0665563, a4304dfa, ba3280c9 and
# cc28ea94 (rows 4-7) had just been handed off, so their Priority flips
# from the default "low" to "ht", and the Latest Handoff Datetime is
# refreshed to the moment the handoff report was generated.

$ws1 = $wb.Worksheets.Item("zh-cn")
for ($r = 4; $r -le 7; $r++) {
    $ws1.Range("E" + $r).Value = "ht"
    $ws1.Range("H" + $r).Value = "2016-08-20 06:39:27"
}

$ws2 = $wb.Worksheets.Item("de-de")
for ($r = 4; $r -le 7; $r++) {
    $ws2.Range("E" + $r).Value = "ht"
    $ws2.Range("H" + $r).Value = "2016-08-20 06:39:31"
}
